# -----------------------------------------------------------------------
# Apply commit: "New files, changes to files and pipeline after bringing
# in memory optimized code."
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. Duplicate the current "High Priority break-up" sheet (sheet 5) BEFORE
#    changing its data, so the duplicate preserves the old numbers. This
#    duplicate becomes the new sheet 6, "Major update - High Priority ".
# -----------------------------------------------------------------------
$wsOldBreakUp = $wb.Worksheets.Item(5)
$wsOldBreakUp.Copy([System.Reflection.Missing]::Value, $wsOldBreakUp)
$wsMajorUpdate = $wb.Worksheets.Item(6)
$wsMajorUpdate.Name = "Major update - High Priority "

# -----------------------------------------------------------------------
# 2. Rename sheet 5 and rewrite its contents (new "Interannual update"
#    numbers, with an extra "Trend New" row inserted above "IUCN").
# -----------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "Interannual update - High Pri"

$ws5.Rows.Item(2).Insert()
$ws5.Rows.Item(2).ClearFormats()

$ws5.Range("A2").Value = "Trend New"
$ws5.Range("B2").Value = 70
$ws5.Range("C2").Value = 68
$ws5.Range("D2").Value = 70
$ws5.Range("E2").Value = 76.09999999999999

$ws5.Range("B3").Value = 33
$ws5.Range("C3").Value = 32
$ws5.Range("D3").Value = 22
$ws5.Range("E3").Value = 23.9

# -----------------------------------------------------------------------
# 3. "Trends Status" sheet (sheet 1) data updates.
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C2").Value = 0
$ws1.Range("E2").Value = 0

$ws1.Range("C3").Value = 2
$ws1.Range("E3").Value = 10

$ws1.Range("C4").Value = 16
$ws1.Range("E4").Value = 80

$ws1.Range("C5").Value = 2
$ws1.Range("E5").Value = 10

$ws1.Range("C7").Value = 32

$ws1.Range("B8").Value = 412
$ws1.Range("C8").Value = 360

# -----------------------------------------------------------------------
# 4. "Priority Status" sheet (sheet 3) data updates.
# -----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# -----------------------------------------------------------------------
# 5. "Species qualification" sheet (sheet 4) data updates.
# -----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 412
$ws4.Range("B4").Value = 52
$ws4.Range("C4").Value = 20
